# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the "Conversión del día" summary text on Hoja1!A1 with the new rates.
$newText = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 1.85 = 6736.13 pesos
✅ 6736.13 pesos = 1.85 = 938.81 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@

$wsHoja1.Range("A1").Value = $newText

# Update the "tasas" sheet rate figures.
$wsTasas.Range("N10").Value = 542
$wsTasas.Range("O10").Value = 3650.98
$wsTasas.Range("N12").Value = 3645
$wsTasas.Range("O12").Value = 508
